$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9866521954536438
$ws.Range("B1").Value = 1.383930444717407
$ws.Range("C1").Value = 5.736277103424072
$ws.Range("D1").Value = 1.682805299758911
$ws.Range("E1").Value = 1.033794641494751
